$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.181.34"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.908.27"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'312.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5084"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").Value = "'0.3936"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "'0.09328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("D10").Value = "'1.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'41.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "'6.396"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "'20.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "1.898.66"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'7.313"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "'0.9981"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'0.00001123"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'92.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'0.06585"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'17.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").Value = "'0.9987"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'6.231"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "28.216.77"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'11.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "'2.644"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.47%  "
$ws.Range("D27").Value = "2.125.08"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").Value = "'21.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "'157.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'127.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'1.094"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").Value = "'0.1071"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'5.627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "'3.607"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'9.655"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").Value = "'0.06689"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'0.02414"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "'1.240"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'0.2187"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +7.42%  "
$ws.Range("D41").Value = "'0.6405"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "'0.9986"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'13.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "'0.6014"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'3.705"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").Value = "'2.023"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").Value = "'123.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "'1.182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
